# Update crypto price (D) and volume-change (E) columns with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.251.53"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "'2.842.86"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'360.55"
$ws.Range("E5").Value = "  +5.57%  "
$ws.Range("D6").Value = "'113.63"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").Value = "'0.575"
$ws.Range("E7").Value = "  +3.82%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").Value = "'41.62"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "'3.294.25"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "'2.821.19"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "'52.175.54"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "'7.64"
$ws.Range("E19").Value = "  +9.60%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "'13.56"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").Value = "'70.44"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'269.02"
$ws.Range("E24").Value = "  -3.52%  "
$ws.Range("D25").Value = "'2.84"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "'27.17"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "'10.43"
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").Value = "'53.88"
$ws.Range("E30").Value = "  +6.78%  "
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "'0.0473"
$ws.Range("E32").Value = "  +24.74%  "
$ws.Range("D33").Value = "'34.65"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").Value = "'5.40"
$ws.Range("E35").Value = "  +7.74%  "
$ws.Range("D36").Value = "'0.0847"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'3.27"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "'18.46"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("D41").Value = "'23.85"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "'128.25"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("E44").Value = "  -6.16%  "
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").Value = "'2.120.04"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  +10.74%  "
$ws.Range("E50").Value = "  +5.35%  "
$ws.Range("D51").Value = "'62.07"
$ws.Range("E51").Value = "  +2.99%  "
